$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell as literal text (leading apostrophe forces text
# interpretation so numeric-looking strings like "233.40" keep their exact
# formatting instead of being parsed into a Double), then reset the style
# back to Normal so no stray quote-prefix style/number-format is left behind.

$ws.Range("D2").Value = "'37.921.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.27%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.085.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'233.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.37%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.36%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'59.37"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.61%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E9").Value = "'  +2.43%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.46%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.45%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.87%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.777"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.04%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.36%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.120.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.57%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'37.819.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'71.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +3.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'228.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.55%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'172.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.88%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.138"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.97%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.92%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.82%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Hedera"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.0632"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.38%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'InternetComputer(DFINITY)"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'4.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.15%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.45%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.87%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.08%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.34%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0984"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.13%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'99.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.25%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +2.26%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.32%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'16.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +8.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.448.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.26%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.21%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.46%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.85%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.59%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.276.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.12%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'47.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.46%  "
$ws.Range("E51").Style = "Normal"
